# Generate Report for Handoff
# Update the localization-status workbook to reflect that
# f02b04e1-c7d8-443c-95f0-c0079bb6c5f1.md is now ready for handoff
# (a new handoff was just generated for both the zh-cn and de-de targets).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-10 20:52:44"

# --- de-de sheet ----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-10 20:52:49"
